# Add a new demo log entry: "Precomputed Atmosphere Rendering"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the formatting of the last existing row (row 12) into the new
# row 13 so the date column keeps its date number format and the notes
# column keeps its wrap-text formatting.
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D13").Insert(-4121)   # xlShiftDown
$ws.Application.CutCopyMode = $false

# Fill in the new row's data.
$ws.Range("A13").Value = 44012
$ws.Range("B13").Value = "充分理解在屏幕空间计算纹理LOD的原理：需要结合微积分教材第1066页推导雅可比的过程来理解，尤其是斜边向量的定义，这是算法的关键解。"
$ws.Range("C13").Value = "需要认真的思考算法细节，尤其是数学层面的东西"
$ws.Range("D13").ClearContents()

$ws.Rows.Item(13).RowHeight = 27.6

# Update the saved selection to match.
$ws.Range("C13").Select()
